$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 6).Value = 1.72
$ws.Cells.Item(2, 7).Value = 1.75
$ws.Cells.Item(2, 9).Value = 4.8
$ws.Cells.Item(2, 10).Value = 4.6
$ws.Cells.Item(2, 11).Value = 4.8
$ws.Cells.Item(2, 12).Value = 1.25
$ws.Cells.Item(2, 14).Value = 7
$ws.Cells.Item(2, 15).Value = 1.15
$ws.Cells.Item(2, 16).Value = 3
$ws.Cells.Item(2, 17).Value = 1.46
$ws.Cells.Item(2, 18).Value = 1.82
$ws.Cells.Item(2, 19).Value = 2.16
$ws.Cells.Item(2, 20).Value = 1.53
$ws.Cells.Item(2, 21).Value = 2.76
$ws.Cells.Item(2, 23).Value = 2.34
$ws.Cells.Item(2, 24).Value = 34
$ws.Cells.Item(2, 28).Value = 16
$ws.Cells.Item(2, 30).Value = 19
$ws.Cells.Item(2, 31).Value = 46
$ws.Cells.Item(2, 32).Value = 15
$ws.Cells.Item(2, 35).Value = 42
$ws.Cells.Item(2, 36).Value = 19.5
$ws.Cells.Item(2, 37).Value = 14.5
$ws.Cells.Item(2, 40).Value = 6.2
$ws.Cells.Item(2, 41).Value = 30
$ws.Cells.Item(3, 12).Value = 1.32
$ws.Cells.Item(3, 18).Value = 1.57
$ws.Cells.Item(3, 19).Value = 2.62
$ws.Cells.Item(3, 20).Value = 1.63
$ws.Cells.Item(3, 35).Value = 44
$ws.Cells.Item(3, 40).Value = 9.199999999999999
$ws.Cells.Item(4, 6).Value = 1.77
$ws.Cells.Item(4, 7).Value = 1.92
$ws.Cells.Item(4, 9).Value = 4.6
$ws.Cells.Item(4, 11).Value = 4.9
$ws.Cells.Item(4, 12).Value = 1.26
$ws.Cells.Item(4, 14).Value = 4.9
$ws.Cells.Item(4, 16).Value = 2.38
$ws.Cells.Item(4, 18).Value = 1.53
$ws.Cells.Item(4, 22).Value = 1.27
$ws.Cells.Item(4, 23).Value = 2.08
$ws.Cells.Item(4, 25).Value = 24
$ws.Cells.Item(4, 28).Value = 14.5
$ws.Cells.Item(4, 29).Value = 11.5
$ws.Cells.Item(5, 7).Value = 2.92
$ws.Cells.Item(5, 8).Value = 3
$ws.Cells.Item(5, 9).Value = 4.3
$ws.Cells.Item(5, 17).Value = 2.12
$ws.Cells.Item(5, 20).Value = 1.05
$ws.Cells.Item(5, 22).Value = 1.3
$ws.Cells.Item(5, 23).Value = 1.52
$ws.Cells.Item(6, 6).Value = 1.5
$ws.Cells.Item(6, 7).Value = 1.8
$ws.Cells.Item(6, 8).Value = 5.5
$ws.Cells.Item(6, 9).Value = 16
$ws.Cells.Item(6, 10).Value = 3.75
$ws.Cells.Item(6, 11).Value = 6.8
$ws.Cells.Item(6, 12).Value = 1.35
$ws.Cells.Item(6, 14).Value = 2.74
$ws.Cells.Item(6, 16).Value = 1.8
$ws.Cells.Item(6, 17).Value = 1.87
$ws.Cells.Item(6, 18).Value = 1.25
$ws.Cells.Item(6, 20).Value = 1.05
$ws.Cells.Item(6, 22).Value = 1.07
$ws.Cells.Item(6, 23).Value = 2.24
$ws.Cells.Item(8, 6).Value = 2.16
$ws.Cells.Item(8, 8).Value = 1.36
$ws.Cells.Item(8, 9).Value = 4.3
$ws.Cells.Item(8, 10).Value = 3
$ws.Cells.Item(8, 13).Value = 1.07
$ws.Cells.Item(8, 14).Value = 1.64
$ws.Cells.Item(8, 15).Value = 1.07
$ws.Cells.Item(8, 16).Value = 1.64
$ws.Cells.Item(8, 20).Value = 1.05
$ws.Cells.Item(8, 22).Value = 1.3
$ws.Cells.Item(9, 6).Value = 2.3
$ws.Cells.Item(9, 7).Value = 2.32
$ws.Cells.Item(9, 13).Value = 1.09
$ws.Cells.Item(9, 20).Value = 1.9
$ws.Cells.Item(9, 21).Value = 2.08
$ws.Cells.Item(9, 22).Value = 1.38
$ws.Cells.Item(9, 24).Value = 12
$ws.Cells.Item(10, 8).Value = 6.8
$ws.Cells.Item(10, 9).Value = 7
$ws.Cells.Item(10, 13).Value = 0
$ws.Cells.Item(10, 19).Value = 2.7
$ws.Cells.Item(10, 20).Value = 1.77
$ws.Cells.Item(10, 21).Value = 2.22
$ws.Cells.Item(10, 23).Value = 2.76
$ws.Cells.Item(11, 13).Value = 1.03
$ws.Cells.Item(11, 17).Value = 1.54
$ws.Cells.Item(11, 19).Value = 2.3
$ws.Cells.Item(12, 7).Value = 2.32
$ws.Cells.Item(12, 10).Value = 3.75
$ws.Cells.Item(12, 14).Value = 4
$ws.Cells.Item(12, 15).Value = 1.31
$ws.Cells.Item(12, 21).Value = 2.24
$ws.Cells.Item(12, 31).Value = 38
$ws.Cells.Item(12, 32).Value = 14.5
$ws.Cells.Item(12, 34).Value = 17
$ws.Cells.Item(12, 37).Value = 23
$ws.Cells.Item(12, 40).Value = 17
$ws.Cells.Item(13, 6).Value = 9.199999999999999
$ws.Cells.Item(13, 7).Value = 9.4
$ws.Cells.Item(13, 14).Value = 5.1
$ws.Cells.Item(13, 23).Value = 1.11
$ws.Cells.Item(13, 26).Value = 8.4
